# fix: typos on menu
# Update the "Rekesalat" and "Vegetarsalat" descriptions: switch the verb
# form from "Serveres ... hjemmelaget" to "Serverast ... heimelaga".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Reker, salat, agurk, mais, fersken, ananas. Serverast med heimelaga dressing og brød. 1,2,3,5,8"
$ws.Range("C3").Value = "Salat, agurk, paprika, mais, ananas, fersken, kikerter, bønner, soltørka tomat. Serverast med heimelaga dressing og brød. 1,2,8"

# Move the selection to C3 (the cell that was last edited), matching the
# saved cursor position in the workbook.
$ws.Range("C3").Select()
